$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1859.4828
$ws.Range("I17").Value = 1307.6923
$ws.Range("J17").Value = 2307.8125
$ws.Range("K17").Value = 3923.0769
$ws.Range("L17").Value = 6923.4375
$ws.Range("M17").Value = -3755.0769
$ws.Range("N17").Value = -7259.4375

# Row 70
$ws.Range("H70").Value = 4678.1333
$ws.Range("I70").Value = 2195.5715
$ws.Range("J70").Value = 6850.375
$ws.Range("K70").Value = 6586.7145
$ws.Range("L70").Value = 20551.125
$ws.Range("M70").Value = -6316.7145
$ws.Range("N70").Value = -21091.125

# Row 73
$ws.Range("H73").Value = 4678.1333
$ws.Range("I73").Value = 2195.5715
$ws.Range("J73").Value = 6850.375
$ws.Range("K73").Value = 6586.7145
$ws.Range("L73").Value = 20551.125
$ws.Range("M73").Value = -5650.7145
$ws.Range("N73").Value = -22423.125

# Row 100
$ws.Range("H100").Value = 3167.5715
$ws.Range("I100").Value = 1149.5555
$ws.Range("J100").Value = 6800
$ws.Range("K100").Value = 1149.5555
$ws.Range("L100").Value = 6800
$ws.Range("M100").Value = -608.5554999999999
$ws.Range("N100").Value = -7882

# Row 107
$ws.Range("H107").Value = 216.25
$ws.Range("I107").Value = 122.166664
$ws.Range("K107").Value = 122.166664
$ws.Range("M107").Value = 1797.833336

# Row 116
$ws.Range("H116").Value = 3848.1
$ws.Range("I116").Value = 3593.8
$ws.Range("J116").Value = 4102.4
$ws.Range("K116").Value = 3593.8
$ws.Range("L116").Value = 4102.4
$ws.Range("M116").Value = -151.8000000000002
$ws.Range("N116").Value = -10986.4

# Row 137
$ws.Range("H137").Value = 2606.5151
$ws.Range("I137").Value = 1945.6875
$ws.Range("J137").Value = 3228.4707
$ws.Range("K137").Value = 5837.0625
$ws.Range("L137").Value = 9685.4121
$ws.Range("M137").Value = -3287.0625
$ws.Range("N137").Value = -14785.4121

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3178.5454
$ws.Range("I45").Value = 2245.5
$ws.Range("K45").Value = 2245.5
$ws.Range("M45").Value = -1868.5

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 4237.2
$ws.Range("I99").Value = 3770.4285
$ws.Range("K99").Value = 3770.4285
$ws.Range("M99").Value = -2272.4285

# Row 107
$ws.Range("H107").Value = 4363.9165
$ws.Range("I107").Value = 902.8461
$ws.Range("K107").Value = 902.8461
$ws.Range("M107").Value = 1017.1539

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Row 134
$ws.Range("H134").Value = 3528.9524
$ws.Range("I134").Value = 2807.3333
$ws.Range("K134").Value = 8421.999899999999
$ws.Range("M134").Value = -5886.999899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 999.4
$ws.Range("I16").Value = 999.5
$ws.Range("J16").Value = 999
$ws.Range("K16").Value = 999.5
$ws.Range("L16").Value = 999
$ws.Range("M16").Value = -712.5
$ws.Range("N16").Value = -1573

# Row 22
$ws.Range("H22").Value = 4399.8
$ws.Range("I22").Value = 3999.6667
$ws.Range("K22").Value = 3999.6667
$ws.Range("M22").Value = -3649.6667

# Row 51
$ws.Range("H51").Value = 27225
$ws.Range("I51").Value = 8500
$ws.Range("J51").Value = 39708.332
$ws.Range("K51").Value = 8500
$ws.Range("L51").Value = 39708.332
$ws.Range("M51").Value = -7764
$ws.Range("N51").Value = -41180.332

# Row 60
$ws.Range("H60").Value = 16833.334
$ws.Range("I60").Value = 2000
$ws.Range("J60").Value = 24250
$ws.Range("K60").Value = 2000
$ws.Range("L60").Value = 24250
$ws.Range("M60").Value = -1489
$ws.Range("N60").Value = -25272

# Row 61
$ws.Range("H61").Value = 27225
$ws.Range("I61").Value = 8500
$ws.Range("J61").Value = 39708.332
$ws.Range("K61").Value = 8500
$ws.Range("L61").Value = 39708.332
$ws.Range("M61").Value = -8152
$ws.Range("N61").Value = -40404.332

# Row 113
$ws.Range("H113").Value = 999.4
$ws.Range("I113").Value = 999.5
$ws.Range("J113").Value = 999
$ws.Range("K113").Value = 999.5
$ws.Range("L113").Value = 999
$ws.Range("M113").Value = 1170.5
$ws.Range("N113").Value = -5339

# Row 122
$ws.Range("H122").Value = 1844.7273
$ws.Range("I122").Value = 1831.2
$ws.Range("J122").Value = 1980
$ws.Range("K122").Value = 5493.6
$ws.Range("L122").Value = 5940
$ws.Range("M122").Value = -3043.6
$ws.Range("N122").Value = -10840

# Row 141
$ws.Range("H141").Value = 87399
$ws.Range("J141").Value = 87399
$ws.Range("L141").Value = 87399
$ws.Range("N141").Value = -97759

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 1519
$ws.Range("J121").Value = 2523.375
$ws.Range("L121").Value = 7570.125
$ws.Range("N121").Value = -10190.125

# Row 131
$ws.Range("H131").Value = 2383.375
$ws.Range("J131").Value = 2612.375
$ws.Range("L131").Value = 7837.125
$ws.Range("N131").Value = -17917.125

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 205336.6
$ws.Range("I132").Value = 255420.75
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 766262.25
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -763732.25
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 6197.391
$ws.Range("J100").Value = 8059
$ws.Range("L100").Value = 8059
$ws.Range("N100").Value = -9141

# Row 133
$ws.Range("H133").Value = 1979899
$ws.Range("J133").Value = 1979899
$ws.Range("L133").Value = 1979899
$ws.Range("N133").Value = -1984959

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 658.75
$ws.Range("I100").Value = 642.5
$ws.Range("J100").Value = 675
$ws.Range("K100").Value = 1285
$ws.Range("L100").Value = 1350
$ws.Range("M100").Value = -744
$ws.Range("N100").Value = -2432

# Row 127
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()

# Row 132
$ws.Range("H132").Value = 2646
$ws.Range("I132").Value = 2646
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7938
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5408
$ws.Range("N132").ClearContents()

